$wb = $excel.ActiveWorkbook
$enSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $enSheet)
$ws.Name = "es"

# Header row (same as en sheet)
$ws.Cells.Item(1,1).Value = 'Key'
$ws.Cells.Item(1,2).Value = 'Value'
$ws.Cells.Item(1,3).Value = 'VoiceDuration'
$ws.Cells.Item(1,4).Value = 'MaxChars'

$ws.Cells.Item(2,1).Value = 'welcome'
$ws.Cells.Item(2,2).Value = '¡Bienvenido!'
$ws.Cells.Item(3,1).Value = 'test'
$ws.Cells.Item(3,2).Value = 'Prueba 123 123'
$ws.Cells.Item(4,1).Value = 'title'
$ws.Cells.Item(4,2).Value = 'Pengu\n<size=30>y el</size>\nEncuentros improbables'
$ws.Cells.Item(5,1).Value = 'play'
$ws.Cells.Item(5,2).Value = 'JUGAR'
$ws.Cells.Item(6,1).Value = 'credits'
$ws.Cells.Item(6,2).Value = 'CRÉDITOS'
$ws.Cells.Item(7,1).Value = 'credits_desc'
$ws.Cells.Item(7,2).Value = 'Escrito por: David Dionisio\r\nMúsica de: Kevin Macleod'
$ws.Cells.Item(8,1).Value = 'options'
$ws.Cells.Item(8,2).Value = 'OPCIONES'
$ws.Cells.Item(9,1).Value = 'sound'
$ws.Cells.Item(9,2).Value = 'SONIDO'
$ws.Cells.Item(10,1).Value = 'music'
$ws.Cells.Item(10,2).Value = 'MÚSICA'
$ws.Cells.Item(11,1).Value = 'speech'
$ws.Cells.Item(11,2).Value = 'DISCURSO'
$ws.Cells.Item(12,1).Value = 'on'
$ws.Cells.Item(12,2).Value = 'PRENDIDO'
$ws.Cells.Item(13,1).Value = 'off'
$ws.Cells.Item(13,2).Value = 'APAGADO'
$ws.Cells.Item(14,1).Value = 'close'
$ws.Cells.Item(14,2).Value = 'CERCA'
$ws.Cells.Item(15,1).Value = 'submit'
$ws.Cells.Item(15,2).Value = 'ENVIAR'
$ws.Cells.Item(16,1).Value = 'victory'
$ws.Cells.Item(16,2).Value = 'VICTORIA'
$ws.Cells.Item(17,1).Value = 'attack'
$ws.Cells.Item(17,2).Value = 'ATAQUE'
$ws.Cells.Item(18,1).Value = 'defend'
$ws.Cells.Item(18,2).Value = 'DEFENDER'
$ws.Cells.Item(19,1).Value = 'rounds'
$ws.Cells.Item(19,2).Value = 'RONDAS'
$ws.Cells.Item(20,1).Value = 'complete'
$ws.Cells.Item(20,2).Value = 'COMPLETAR'
$ws.Cells.Item(21,1).Value = 'revivePenalty'
$ws.Cells.Item(21,2).Value = 'REVIVIR PENALIZACIÓN'
$ws.Cells.Item(22,1).Value = 'total'
$ws.Cells.Item(22,2).Value = 'TOTAL'
$ws.Cells.Item(23,1).Value = 'level_1_intro_1'
$ws.Cells.Item(23,2).Value = 'Un abismo impide el viaje de Pengu.'
$ws.Cells.Item(24,1).Value = 'level_1_intro_2'
$ws.Cells.Item(24,2).Value = '¡Ayuda a Pengu a salir determinando la distancia de la brecha usando fracciones!'
$ws.Cells.Item(25,1).Value = 'level_1_info_1'
$ws.Cells.Item(25,2).Value = 'Para sumar estas dos fracciones, debes hacer iguales los denominadores.'
$ws.Cells.Item(26,1).Value = 'level_1_info_2_a'
$ws.Cells.Item(26,2).Value = 'Una forma de hacerlo es multiplicando los denominadores juntos para hacerlos iguales.'
$ws.Cells.Item(27,1).Value = 'level_1_info_2_b'
$ws.Cells.Item(28,1).Value = 'level_1_info_2_c'
$ws.Cells.Item(28,2).Value = 'Después de eso, puedes agregar ambas fracciones correctamente.'
$ws.Cells.Item(29,1).Value = 'level_1_info_3'
$ws.Cells.Item(29,2).Value = '¡Ahora te toca a ti! Utilice el multiplicador para hacer que ambos denominadores sean iguales.'
$ws.Cells.Item(30,1).Value = 'level_1_info_4'
$ws.Cells.Item(30,2).Value = '¿Observa cómo ambas fracciones tienen ahora los mismos tamaños de unidad? Ahora se pueden añadir correctamente.'
$ws.Cells.Item(31,1).Value = 'level_1_info_5'
$ws.Cells.Item(31,2).Value = 'Sigue adelante y escribe la respuesta correcta pulsando en cualquiera de las ranuras.'
$ws.Cells.Item(32,1).Value = 'success'
$ws.Cells.Item(32,2).Value = '¡ÉXITO!'
$ws.Cells.Item(33,1).Value = 'card_drag_instruct'
$ws.Cells.Item(33,2).Value = 'Arrastra una tarjeta a una ranura vacía.'
$ws.Cells.Item(34,1).Value = 'level_2_intro_1'
$ws.Cells.Item(34,2).Value = 'Una vez más, un obstáculo bloquea el camino de Pengu. Esta vez, una foca elefante.'
$ws.Cells.Item(35,1).Value = 'level_2_intro_2'
$ws.Cells.Item(35,2).Value = '¡Golpea el sello con suficientes fracciones para sacarlo del camino!'
$ws.Cells.Item(36,1).Value = 'level_3_intro_1'
$ws.Cells.Item(36,2).Value = '¡Cuidado! ¡Un yeti se interviene en camino!'
$ws.Cells.Item(37,1).Value = 'level_3_intro_2'
$ws.Cells.Item(37,2).Value = '¡Debemos soportar el viento helado para continuar!'
$ws.Cells.Item(38,1).Value = 'level_3_info_1'
$ws.Cells.Item(38,2).Value = 'Para tener éxito, debes restar el ataque fraccional hasta que llegue por debajo de cero.'
$ws.Cells.Item(39,1).Value = 'mixedNumber'
$ws.Cells.Item(39,2).Value = 'Número Mixto'
$ws.Cells.Item(40,1).Value = 'improperFraction'
$ws.Cells.Item(40,2).Value = 'Fracción Incorrecta'
$ws.Cells.Item(41,1).Value = 'level_5_intro_1'
$ws.Cells.Item(41,2).Value = 'Después de una onerosa batalla, Pengu debe respirar aire fresco.'
$ws.Cells.Item(42,1).Value = 'level_5_intro_2'
$ws.Cells.Item(42,2).Value = '¡Ayuda a Pengu a nadar hacia la tierra!'
$ws.Cells.Item(43,1).Value = 'level_5_info_1_a'
$ws.Cells.Item(43,2).Value = 'Como puede ver, hay un número mixto en esta operación.'
$ws.Cells.Item(44,1).Value = 'level_5_info_1_b'
$ws.Cells.Item(44,2).Value = 'Un número mixto se compone de un número entero y una fracción.'
$ws.Cells.Item(45,1).Value = 'level_5_info_2_a'
$ws.Cells.Item(45,2).Value = 'Para convertir un número mixto en una fracción incorrecta: multiplique el número entero por el denominador de la fracción.'
$ws.Cells.Item(46,1).Value = 'level_5_info_2_b'
$ws.Cells.Item(46,2).Value = 'Después, suma el resultado al numerador.'
$ws.Cells.Item(47,1).Value = 'level_5_info_3_a'
$ws.Cells.Item(47,2).Value = 'Puede arrastrar el número entero hacia la fracción, o viceversa, para convertirlo.'
$ws.Cells.Item(48,1).Value = 'level_5_info_3_b'
$ws.Cells.Item(48,2).Value = '¡Utiliza esta técnica para ayudarte con operaciones complicadas!'
$ws.Cells.Item(49,1).Value = 'level_7_intro_1'
$ws.Cells.Item(49,2).Value = '¡Una roca está bloqueando el camino de Pengu!'
$ws.Cells.Item(50,1).Value = 'level_7_intro_2'
$ws.Cells.Item(50,2).Value = '¡Es hora de desatar los poderes fraccionarios más potentes de Pengu!'
$ws.Cells.Item(51,1).Value = 'game_complete'
$ws.Cells.Item(51,2).Value = '¡En hora buena! ¡Has descubierto el tesoro!'
$ws.Cells.Item(52,1).Value = 'game_complete_desc'
$ws.Cells.Item(52,2).Value = '¡Gracias por jugar!'
$ws.Cells.Item(53,1).Value = 'finish'
$ws.Cells.Item(53,2).Value = 'ACABADO'
$ws.Cells.Item(54,1).Value = 'total_score'
$ws.Cells.Item(54,2).Value = 'PUNTOS TOTALES:'

# Numeric columns preserved from en sheet
$ws.Cells.Item(2,3).Value = 1.5
$ws.Cells.Item(2,4).Value = 50
$ws.Cells.Item(51,3).Value = 4
$ws.Cells.Item(52,3).Value = 3